$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Message bodies for column F (column 6), reused across several rows
$msg1 = "DEMO_POST`n🔥 **DEMO DEALS OF THE DAY** 🔥`n💥 Price Crash Store`n⚡️ Up to 5% off`n👉🏻 amzaff.in/l3swo0g`n🌟 Kid's Carnival`n📣 Sale live till 28th July`n👉🏻 amzaff.in/jrtPYsT`n🧸 Toy's Fiesta`n⚡️ Up to 70% off`n👉🏻 amzaff.in/pXpqAHe`n🪴 Solar Garden & Outdoors`n⚡️ Up to 70% off`n👉🏻 amzaff.in/092HVM2`n🛍 End Of Season Sale`n⚡️ Up to 70% off`n👉🏻 amzaff.in/ukuS1wj`n🏠 Home Shopping Spree`n📣 Sale live from 26th-30th July`n👉🏻 amzaff.in/dvgN1JH"
$msg2 = "DEMO_POST`n💧 **DEMO Essentials for Skin, Hair & Fragrance** 💧`n💥 Maximise earnings with Beauty commissions – Now increased to 10%`n🌿 Mamaearth Rice Oil-Free Face Moisturizer for Oily Skin`n⚡️ 80g @25% + 5% Off – ₹22`n👉🏻 amzaff.in/kyKGkVq`n🧼 Cetaphil Brightness Reveal Creamy Cleanser`n⚡️ 100g @25% Off – ₹599`n👉🏻 amzaff.in/2QkCAT6`n💆‍♀️ Herbal Essences bio:renew Argan Oil of Morocco Shampoo`n⚡️ 400ml @60% Off – ₹260`n👉🏻 amzaff.in/WPdHenG`n🧴 Be Bodywise 6% AHA BHA Underarm Roll On Deodorant`n⚡️ 50ml | Alcohol & Aluminum Free – ₹399`n👉🏻 amzaff.in/3YBXpxC`n🕺 Park Avenue Voyage Signature Collection Perfume for Men`n⚡️ 120ml @57% Off – ₹171`n👉🏻 amzaff.in/IQsEdXX"
$msg3 = "DEMO_POST`n💻 **DEMO Level Gaming Laptops**`n⚡️ Up to 45% off`n👉🏻 amzaff.in/FeVABNi`n🔥 Gaming Laptops Under ₹80,000 🔥`n💻 ASUS TUF A15 GAMING`n⚡️ ~~₹83,990~~ | ₹63,990`n⚡️ Effective price ₹56,490`n👉🏻 amzaff.in/HNJJ4b3`n💻 HP VICTUS GAMING`n⚡️ ~~₹99,382~~ | ₹82,990`n⚡️ Effective price ₹77,490`n👉🏻 amzaff.in/rIYTx8U`n💻 LENOVO LOQ GAMING`n⚡️ ~~₹1,12,990~~ | ₹87,190`n⚡️ Effective price ₹77,190`n👉🏻 amzaff.in/Q2oXx7g`n💻 ACER NITRO V GAMING`n⚡️ ~~₹89,999~~ | ₹70,990`n⚡️ Effective price ₹67,490`n👉🏻 amzaff.in/Es6mtU0`n💻 DELL G15 GAMING`n⚡️ ~~₹1,06,331~~ | ₹77,490`n⚡️ Effective price ₹68,990`n👉🏻 amzaff.in/MMEYXBc`n💻 ASUS CREATOR SERIES`n⚡️ ~~₹85,990~~ | ₹72,990`n⚡️ Effective price ₹66,490`n👉🏻 amzaff.in/R4IV7C8"

# New scheduler rows appended below the existing 4 rows (rows 5-13).
# Column C holds a plain "YYYY-MM-DD" text label (not a real date), so each
# cell is briefly formatted as Text before the assignment to stop Excel from
# auto-converting the literal string into a date serial number; the style is
# then reset back to Normal so formatting matches the rest of the sheet.

$ws.Cells.Item(5, 1).Value = 1
$ws.Cells.Item(5, 2).Value = "Kid's Carnival"
$ws.Cells.Item(5, 3).NumberFormat = "@"
$ws.Cells.Item(5, 3).Value = "2025-08-06"
$ws.Cells.Item(5, 3).Style = "Normal"
$ws.Cells.Item(5, 4).Value = "14:30:00"
$ws.Cells.Item(5, 5).Value = "✅ Scheduled"
$ws.Cells.Item(5, 6).Value = $msg1

$ws.Cells.Item(6, 1).Value = 2
$ws.Cells.Item(6, 2).Value = "Daily Essentials"
$ws.Cells.Item(6, 3).NumberFormat = "@"
$ws.Cells.Item(6, 3).Value = "2025-08-06"
$ws.Cells.Item(6, 3).Style = "Normal"
$ws.Cells.Item(6, 4).Value = "15:15:00"
$ws.Cells.Item(6, 5).Value = "✅ Scheduled"
$ws.Cells.Item(6, 6).Value = $msg2

$ws.Cells.Item(7, 1).Value = 3
$ws.Cells.Item(7, 2).Value = "Laptops"
$ws.Cells.Item(7, 3).NumberFormat = "@"
$ws.Cells.Item(7, 3).Value = "2025-08-06"
$ws.Cells.Item(7, 3).Style = "Normal"
$ws.Cells.Item(7, 4).Value = "15:30:00"
$ws.Cells.Item(7, 5).Value = "✅ Scheduled"
$ws.Cells.Item(7, 6).Value = $msg3

$ws.Cells.Item(8, 1).Value = 1
$ws.Cells.Item(8, 2).Value = "Kid's Carnival"
$ws.Cells.Item(8, 3).NumberFormat = "@"
$ws.Cells.Item(8, 3).Value = "2025-08-06"
$ws.Cells.Item(8, 3).Style = "Normal"
$ws.Cells.Item(8, 4).Value = "14:30:00"
$ws.Cells.Item(8, 5).Value = "✅ Scheduled"
$ws.Cells.Item(8, 6).Value = $msg1

$ws.Cells.Item(9, 1).Value = 2
$ws.Cells.Item(9, 2).Value = "Daily Essentials"
$ws.Cells.Item(9, 3).NumberFormat = "@"
$ws.Cells.Item(9, 3).Value = "2025-08-06"
$ws.Cells.Item(9, 3).Style = "Normal"
$ws.Cells.Item(9, 4).Value = "15:00:00"
$ws.Cells.Item(9, 5).Value = "✅ Scheduled"
$ws.Cells.Item(9, 6).Value = $msg2

$ws.Cells.Item(10, 1).Value = 3
$ws.Cells.Item(10, 2).Value = "Laptops"
$ws.Cells.Item(10, 3).NumberFormat = "@"
$ws.Cells.Item(10, 3).Value = "2025-08-06"
$ws.Cells.Item(10, 3).Style = "Normal"
$ws.Cells.Item(10, 4).Value = "15:30:00"
$ws.Cells.Item(10, 5).Value = "✅ Scheduled"
$ws.Cells.Item(10, 6).Value = $msg3

$ws.Cells.Item(11, 1).Value = 1
$ws.Cells.Item(11, 2).Value = "Kid's Carnival"
$ws.Cells.Item(11, 3).NumberFormat = "@"
$ws.Cells.Item(11, 3).Value = "2025-08-05"
$ws.Cells.Item(11, 3).Style = "Normal"
$ws.Cells.Item(11, 4).Value = "15:30:00"
$ws.Cells.Item(11, 5).Value = "✅ Scheduled"
$ws.Cells.Item(11, 6).Value = $msg1

$ws.Cells.Item(12, 1).Value = 2
$ws.Cells.Item(12, 2).Value = "Daily Essentials"
$ws.Cells.Item(12, 3).NumberFormat = "@"
$ws.Cells.Item(12, 3).Value = "2025-08-05"
$ws.Cells.Item(12, 3).Style = "Normal"
$ws.Cells.Item(12, 4).Value = "16:05:00"
$ws.Cells.Item(12, 5).Value = "✅ Scheduled"
$ws.Cells.Item(12, 6).Value = $msg2

$ws.Cells.Item(13, 1).Value = 3
$ws.Cells.Item(13, 2).Value = "Laptops"
$ws.Cells.Item(13, 3).NumberFormat = "@"
$ws.Cells.Item(13, 3).Value = "2025-08-05"
$ws.Cells.Item(13, 3).Style = "Normal"
$ws.Cells.Item(13, 4).Value = "16:30:00"
$ws.Cells.Item(13, 5).Value = "✅ Scheduled"
$ws.Cells.Item(13, 6).Value = $msg3
